$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the unit price in I49 (4.5 -> 2.25)
$ws.Range("I49").Value = 2.25

# Expand the print area from A1:G59 to A1:L59
$ws.PageSetup.PrintArea = "`$A`$1:`$L`$59"

# Reduce the print scale from 72% to 69%
$ws.PageSetup.Zoom = 69

# Move the active cell / selection on the frozen (bottom-right) pane to I49
$ws.Activate()
$ws.Range("I49").Select()
